$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: column D label changes from "Valor (USD$)" to "Valor ($COP)"
$ws.Range("D1").Value = 'Valor ($COP)'

# --- Row 2 ---
# Dates in this sheet are stored as literal text (not real Excel dates), so
# we apply the sheet's existing date number-format to each date cell first
# (this matches/reuses the style the cell already carried) and then assign
# the value with a leading quote so it is kept as literal text rather than
# being auto-parsed into a date serial number.
$ws.Cells.Item(2,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(2,1).Value = "'01/09/25"
$ws.Cells.Item(2,2).Value = 'Compra concentrado peces 1lb'
$ws.Cells.Item(2,3).Value = 'Egreso - Operacional'
$ws.Cells.Item(2,4).Value = 15000
$ws.Cells.Item(2,5).Value = 'Lina'

# --- Row 3 ---
$ws.Cells.Item(3,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(3,1).Value = "'01/09/25"
$ws.Cells.Item(3,2).Value = 'Venta pescado 1kg'
$ws.Cells.Item(3,3).Value = 'Ingreso - Venta'
$ws.Cells.Item(3,4).Value = 40000
$ws.Cells.Item(3,5).Value = 'Lina'

# --- Row 4 ---
$ws.Cells.Item(4,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(4,1).Value = "'02/09/25"
$ws.Cells.Item(4,2).Value = 'Pago aportes Amparo Cano septiembre'
$ws.Cells.Item(4,3).Value = 'Ingreso - Aporte'
$ws.Cells.Item(4,4).Value = 50000
$ws.Cells.Item(4,5).Value = 'Lina'

# --- Row 5 ---
$ws.Cells.Item(5,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(5,1).Value = "'02/09/25"
$ws.Cells.Item(5,2).Value = 'Pago aportes Luz Mary Septiembre'
$ws.Cells.Item(5,3).Value = 'Ingreso - Aporte'
$ws.Cells.Item(5,4).Value = 50000
$ws.Cells.Item(5,5).Value = 'Lina'

# --- Row 6 ---
$ws.Cells.Item(6,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(6,1).Value = "'02/09/25"
$ws.Cells.Item(6,2).Value = 'Pago aportes Sergio Hernán Septiembre'
$ws.Cells.Item(6,3).Value = 'Ingreso - Aporte'
$ws.Cells.Item(6,4).Value = 50000
$ws.Cells.Item(6,5).Value = 'Lina'

# --- Row 7 (new) ---
$ws.Cells.Item(7,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(7,1).Value = "'04/09/25"
$ws.Cells.Item(7,2).Value = 'Venta peces 1lb'
$ws.Cells.Item(7,3).Value = 'Ingreso - Venta'
$ws.Cells.Item(7,4).Value = 20000
$ws.Cells.Item(7,5).Value = 'Lina'

# --- Row 8 (new) ---
$ws.Cells.Item(8,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(8,1).Value = "'05/09/25"
$ws.Cells.Item(8,2).Value = 'Pago servicios públicos'
$ws.Cells.Item(8,3).Value = 'Egreso - Operacional'
$ws.Cells.Item(8,4).Value = 23800
$ws.Cells.Item(8,5).Value = 'Lina'
